# Adds a "2022-Q4" sheet (fund-holdings detail) right after "总计",
# and updates the "总计" summary sheet with the new quarter's totals.
#
# Final sheet order: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 (summary) sheet: insert a new row 2 for "2022-Q4" and push the
#    existing Q3/Q2/Q1 rows down by one; re-create the trailing Q1 row
#    that falls off the bottom (old row 4 -> new row 5).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room: inserting at row 2 shifts old rows 2,3,4 -> 3,4,5
$total.Rows(2).Insert()

# Row 2 lost its formatting on insert; clone it from row 3 (old row 2,
# already shifted down, still carries the correct style/format).
$total.Range("A3:D3").Copy($total.Range("A2:D2"))

# Now write the actual 2022-Q4 totals into row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 14
$total.Range("D2").Value = 1.28

# Rows 3 and 4 already hold the old Q3/Q2 data (shifted down) - update
# their label/index text/values to match the new row meaning is NOT
# needed for values (they keep their own numbers), but the running
# index in column A and date label in column B must be refreshed.
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.25

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.19

# Row 5 is brand new - clone formatting from row 4, then fill values.
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.17

# ---------------------------------------------------------------------
# 2) Insert the brand-new "2022-Q4" sheet right after "总计" (i.e.
#    right before the current first quarter tab, "2022-Q3").
# ---------------------------------------------------------------------
$q3tab = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3tab)
$newSheet.Name = "2022-Q4"

# Re-fetch fresh handles by name (structural ops can invalidate old refs).
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Item("2022-Q4")

# Clone the header row + first two data rows from 2022-Q3 so the new
# sheet inherits identical column styling (bold/centered/bordered header
# in B1:H1, bold/centered/bordered index column A2:A.., plain data
# cells elsewhere).
$q3.Range("A1:H3").Copy($q4.Range("A1:H3"))

# Extend down to 15 rows total (1 header + 14 fund rows) by repeating
# the row-3 formatting pattern.
for ($r = 4; $r -le 15; $r++) {
    $q4.Range("A3:H3").Copy($q4.Range("A" + $r + ":H" + $r))
}

# ---- header row -------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# ---- fund holding rows --------------------------------------------
# columns: A index(n) B code(text) C name(text) D scale(text)
#          E position(text) F ratio(text) G value(text) H rank(n)
$rows = @(
    @(0,  "630008", "华商策略精选混合",             "7.47",  "76.58", "4.75", "0.3548", 2),
    @(1,  "506003", "富国科创板两年定期开放混合",     "13.76", "98.91", "2.18", "0.3000", 10),
    @(2,  "010994", "博时创新经济混合A",             "3.80",  "86.65", "4.56", "0.1733", 5),
    @(3,  "003593", "国泰景气行业灵活配置混合",       "3.51",  "92.21", "3.74", "0.1313", 9),
    @(4,  "013958", "华商鑫选回报一年持有混合A",      "5.33",  "90.81", "2.20", "0.1173", 3),
    @(5,  "010995", "博时创新经济混合C",             "1.71",  "86.65", "4.56", "0.0780", 5),
    @(6,  "016336", "博时卓远成长一年持有期股票A",    "1.63",  "51.15", "4.29", "0.0699", 3),
    @(7,  "013959", "华商鑫选回报一年持有混合C",      "1.22",  "90.81", "2.20", "0.0268", 3),
    @(8,  "016337", "博时卓远成长一年持有期股票C",    "0.46",  "51.15", "4.29", "0.0197", 3),
    @(9,  "008300", "人保量化锐进混合A",             "0.08",  "90.63", "4.41", "0.0035", 3),
    @(10, "008301", "人保量化锐进混合C",             "0.04",  "90.63", "4.41", "0.0018", 3),
    @(11, "166108", "信澳量化多因子混合（LOF）C",     "0.07",  "34.40", "0.66", "0.0005", 8),
    @(12, "166107", "信澳量化多因子混合（LOF）A",     "0.05",  "34.40", "0.66", "0.0003", 8),
    @(13, "004352", "北信瑞丰研究精选股票",           "0.01",  "92.09", "1.86", "0.0002", 4)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $q4.Range("A$r").Value = $data[0]
    $q4.Range("B$r").Value = "'" + $data[1]
    $q4.Range("C$r").Value = $data[2]
    $q4.Range("D$r").Value = "'" + $data[3]
    $q4.Range("E$r").Value = "'" + $data[4]
    $q4.Range("F$r").Value = "'" + $data[5]
    $q4.Range("G$r").Value = "'" + $data[6]
    $q4.Range("H$r").Value = $data[7]

    # Drop the forced-text "quote prefix" formatting artifact picked up
    # above so B, D:G keep the same (unstyled) look as column C / H.
    $q4.Range("B$r").Style = $q4.Range("C$r").Style
    $q4.Range("D$r" + ":G$r").Style = $q4.Range("C$r").Style
}

# Leave the workbook focused back on the summary sheet, as before the edit.
$total.Activate()

